# FMX-1024: Fix clash with ABP id by regeneration.
#
# The source data used to populate column I ("ABP id") of the
# "categories-sous categories" sheet was regenerated, shifting every id
# in rows 145-208 upward by 32 (to avoid clashing with another range of
# ids). The workbook was also re-exported, which appended one more
# generation ("_0_0_0") of the auto-duplicated _xlnm.Print_Area /
# _xlnm.Print_Titles / _xlnm._FilterDatabase defined names, and left the
# active selection on a different cell/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the regenerated ABP ids in column I (rows 145-208) ---
# Rows that hold formulas (=I<row-1>) recompute automatically once the
# seed cell above them changes, so only the "seed" rows need new values.
$ws.Range("I145").Value = 192
$ws.Range("I146").Value = 190
$ws.Range("I147").Value = 193
$ws.Range("I148").Value = 191
$ws.Range("I149").Value = 197
$ws.Range("I150").Value = 195
$ws.Range("I151").Value = 205
$ws.Range("I152").Value = 203
$ws.Range("I153").Value = 181
$ws.Range("I154").Value = 177
$ws.Range("I155").Value = 175
$ws.Range("I156").Value = 173
$ws.Range("I157").Value = 171
$ws.Range("I158").Value = 204
$ws.Range("I159").Value = 202
$ws.Range("I160").Value = 180
$ws.Range("I161").Value = 176
$ws.Range("I162").Value = 174
$ws.Range("I163").Value = 172
$ws.Range("I164").Value = 170
$ws.Range("I165").Value = 167
$ws.Range("I166").Value = 163
$ws.Range("I167").Value = 161
$ws.Range("I168").Value = 159
$ws.Range("I169").Value = 166
$ws.Range("I170").Value = 162
$ws.Range("I171").Value = 160
$ws.Range("I172").Value = 158
$ws.Range("I173").Value = 157
$ws.Range("I179").Value = 156
$ws.Range("I185").Value = 155
$ws.Range("I191").Value = 154
$ws.Range("I197").Value = 153
$ws.Range("I203").Value = 152

# --- Re-generation added one more level of the auto-suffixed builtin
#     defined names (Print_Area / Print_Titles / FilterDatabase) ---
$ws.Names.Add("_xlnm.Print_Area_0_0_0", "='categories-sous categories'!`$A`$1:`$P`$144")
$ws.Names.Add("_xlnm.Print_Titles_0_0_0", "='categories-sous categories'!`$1:`$2")
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0", "='categories-sous categories'!`$A`$2:`$Q`$144")

# --- Move the active selection/scroll position left on the sheet ---
$ws.Activate()
[void]$ws.Range("I204").Select()
$excel.ActiveWindow.ScrollRow = 185
$excel.ActiveWindow.ScrollColumn = 1
